$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-coerced into numbers by Excel (losing trailing zeros / dot-grouping).
$textCells = @("D4","D5","D6","D7","D8","D9","D10","D11","D12","D13","D14","D15","D17","D18","D19","D20","D21","D22","D25","D26","D27","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values.
$ws.Range("D2").Value = '27.513.90'
$ws.Range("E2").Value = '  -1.95%  '
$ws.Range("D3").Value = '1.749.89'
$ws.Range("E3").Value = '  -2.19%  '
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '324.33'
$ws.Range("E5").Value = '  +0.30%  '
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("D7").Value = '0.4476'
$ws.Range("E7").Value = '  +3.84%  '
$ws.Range("D8").Value = '0.3597'
$ws.Range("E8").Value = '  -0.69%  '
$ws.Range("D9").Value = '0.07506'
$ws.Range("E9").Value = '  +0.03%  '
$ws.Range("D10").Value = '41.97'
$ws.Range("E10").Value = '  -6.13%  '
$ws.Range("D11").Value = '1.092'
$ws.Range("E11").Value = '  -1.80%  '
$ws.Range("D12").Value = '1.001'
$ws.Range("D13").Value = '20.61'
$ws.Range("E13").Value = '  -4.93%  '
$ws.Range("D14").Value = '6.021'
$ws.Range("E14").Value = '  -2.01%  '
$ws.Range("D15").Value = '7.123'
$ws.Range("E15").Value = '  -2.77%  '
$ws.Range("D16").Value = '1.755.43'
$ws.Range("E16").Value = '  -0.88%  '
$ws.Range("D17").Value = '93.22'
$ws.Range("E17").Value = '  +0.86%  '
$ws.Range("D18").Value = '0.00001061'
$ws.Range("E18").Value = '  -0.25%  '
$ws.Range("D19").Value = '0.06377'
$ws.Range("E19").Value = '  +0.74%  '
$ws.Range("D20").Value = '1.0000'
$ws.Range("E20").Value = '  +0.03%  '
$ws.Range("D21").Value = '16.76'
$ws.Range("E21").Value = '  -2.66%  '
$ws.Range("D22").Value = '5.844'
$ws.Range("E22").Value = '  -1.83%  '
$ws.Range("D23").Value = '27.561.26'
$ws.Range("E24").Value = '  -1.96%  '
$ws.Range("D25").Value = '2.074'
$ws.Range("E25").Value = '  -2.80%  '
$ws.Range("D26").Value = '161.62'
$ws.Range("E26").Value = '  +1.27%  '
$ws.Range("D27").Value = '20.46'
$ws.Range("E27").Value = '  +0.52%  '
$ws.Range("D28").Value = '1.949.32'
$ws.Range("E28").Value = '  -2.00%  '
$ws.Range("D29").Value = '2.083'
$ws.Range("E29").Value = '  -4.51%  '
$ws.Range("D30").Value = '125.34'
$ws.Range("E30").Value = '  -1.18%  '
$ws.Range("D31").Value = '1.082'
$ws.Range("E31").Value = '  -6.78%  '
$ws.Range("D32").Value = '3.661'
$ws.Range("E32").Value = '  +4.10%  '
$ws.Range("D33").Value = '0.09006'
$ws.Range("E33").Value = '  +0.11%  '
$ws.Range("D34").Value = '5.531'
$ws.Range("E34").Value = '  -3.17%  '
$ws.Range("D35").Value = '11.94'
$ws.Range("E35").Value = '  -5.37%  '
$ws.Range("D36").Value = '0.02294'
$ws.Range("E36").Value = '  -1.15%  '
$ws.Range("D37").Value = '0.06014'
$ws.Range("E37").Value = '  -0.63%  '
$ws.Range("D38").Value = '0.2085'
$ws.Range("E38").Value = '  -1.58%  '
$ws.Range("D39").Value = '0.6337'
$ws.Range("E39").Value = '  -1.70%  '
$ws.Range("D40").Value = '4.941'
$ws.Range("E40").Value = '  -3.03%  '
$ws.Range("E41").Value = '  +1.56%  '
$ws.Range("B42").Value = 'WEMIXTOKEN'
$ws.Range("C42").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D42").Value = '1.379'
$ws.Range("E42").Value = '  -2.63%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '7.751'
$ws.Range("E43").Value = '  -1.42%  '
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Value = '13.14'
$ws.Range("E44").Value = '  -2.79%  '
$ws.Range("B45").Value = 'PancakeSwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D45").Value = '3.712'
$ws.Range("E45").Value = '  +0.27%  '
$ws.Range("D46").Value = '0.5881'
$ws.Range("E46").Value = '  -1.70%  '
$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D47").Value = '122.04'
$ws.Range("E47").Value = '  -1.86%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = '1.955'
$ws.Range("E48").Value = '  -1.41%  '
$ws.Range("B49").Value = 'EOS'
$ws.Range("C49").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D49").Value = '1.144'
$ws.Range("E49").Value = '  -0.60%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '0.06855'
$ws.Range("E50").Value = '  -1.33%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").Value = '71.95'
$ws.Range("E51").Value = '  -3.23%  '
